$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.9180432866803869
$ws.Cells.Item(2, 3).Value = 0.1097131679421821
$ws.Cells.Item(2, 4).Value = 0.06588484050381638
$ws.Cells.Item(2, 6).Value = 3.202136094421334
$ws.Cells.Item(2, 7).Value = 2.454686398660058
$ws.Cells.Item(2, 8).Value = 2.000327881350913
$ws.Cells.Item(2, 10).Value = 0.2631755003609335
$ws.Cells.Item(2, 11).Value = 0.4514623245620726
$ws.Cells.Item(2, 12).Value = 0.2858766676266171
$ws.Cells.Item(2, 13).Value = 0.2582246334812339

$ws.Cells.Item(3, 2).Value = 0.893301933005489
$ws.Cells.Item(3, 3).Value = 0.1081995496284662
$ws.Cells.Item(3, 4).Value = 0.06566430498912013
$ws.Cells.Item(3, 6).Value = 3.198259940160511
$ws.Cells.Item(3, 7).Value = 2.448282934897591
$ws.Cells.Item(3, 8).Value = 2.002215772316163
$ws.Cells.Item(3, 10).Value = 0.2640221064714581
$ws.Cells.Item(3, 11).Value = 0.4273106364189232
$ws.Cells.Item(3, 12).Value = 0.2849207388460755
$ws.Cells.Item(3, 13).Value = 0.2541460332465988

$ws.Cells.Item(4, 2).Value = 0.8786210122522107
$ws.Cells.Item(4, 3).Value = 0.1072459505235557
$ws.Cells.Item(4, 4).Value = 0.06555386845751698
$ws.Cells.Item(4, 6).Value = 3.196985381724474
$ws.Cells.Item(4, 7).Value = 2.44525684978791
$ws.Cells.Item(4, 8).Value = 2.003975794872446
$ws.Cells.Item(4, 10).Value = 0.2645914505129063
$ws.Cells.Item(4, 11).Value = 0.4127471068864139
$ws.Cells.Item(4, 12).Value = 0.2844446862116214
$ws.Cells.Item(4, 13).Value = 0.2517722038863468

$ws.Cells.Item(5, 2).Value = 0.8727671823033631
$ws.Cells.Item(5, 3).Value = 0.1068512432897997
$ws.Cells.Item(5, 4).Value = 0.06551517721165823
$ws.Cells.Item(5, 6).Value = 3.196744162014198
$ws.Cells.Item(5, 7).Value = 2.4442514707915
$ws.Cells.Item(5, 8).Value = 2.004844245172976
$ws.Cells.Item(5, 10).Value = 0.2648359365946398
$ws.Cells.Item(5, 11).Value = 0.4068794218680836
$ws.Cells.Item(5, 12).Value = 0.2842786553075669
$ws.Cells.Item(5, 13).Value = 0.2508377539848219

$ws.Cells.Item(6, 2).Value = 0.8718029480216387
$ws.Cells.Item(6, 3).Value = 0.10678533312824
$ws.Cells.Item(6, 4).Value = 0.06550913479170717
$ws.Cells.Item(6, 6).Value = 3.196720914783256
$ws.Cells.Item(6, 7).Value = 2.444098286594283
$ws.Cells.Item(6, 8).Value = 2.004997588495513
$ws.Cells.Item(6, 10).Value = 0.2648772873077334
$ws.Cells.Item(6, 11).Value = 0.4059091558585379
$ws.Cells.Item(6, 12).Value = 0.2842527773399581
$ws.Cells.Item(6, 13).Value = 0.250684579699108

$ws.Cells.Item(7, 2).Value = 0.8785415435209813
$ws.Cells.Item(7, 3).Value = 0.1072406521122424
$ws.Cells.Item(7, 4).Value = 0.06555332105076417
$ws.Cells.Item(7, 6).Value = 3.196981001970713
$ws.Cells.Item(7, 7).Value = 2.44524236859769
$ws.Cells.Item(7, 8).Value = 2.003986894608232
$ws.Cells.Item(7, 10).Value = 0.2645946972049344
$ws.Cells.Item(7, 11).Value = 0.4126677013195348
$ws.Cells.Item(7, 12).Value = 0.2844423337272275
$ws.Cells.Item(7, 13).Value = 0.2517594682008522

$ws.Cells.Item(8, 2).Value = 0.9094067859004724
$ws.Cells.Item(8, 3).Value = 0.109196292399858
$ws.Cells.Item(8, 4).Value = 0.06580363609417716
$ws.Cells.Item(8, 6).Value = 3.200570284528737
$ws.Cells.Item(8, 7).Value = 2.452290529001488
$ws.Cells.Item(8, 8).Value = 2.000854199859575
$ws.Cells.Item(8, 10).Value = 0.2634571435683917
$ws.Cells.Item(8, 11).Value = 0.4430798207154112
$ws.Cells.Item(8, 12).Value = 0.2855240947277338
$ws.Cells.Item(8, 13).Value = 0.2567913171762228

$ws.Cells.Item(9, 2).Value = 0.9739685960793167
$ws.Cells.Item(9, 3).Value = 0.1128398158241239
$ws.Cells.Item(9, 4).Value = 0.06649130356065314
$ws.Cells.Item(9, 6).Value = 3.216373924029241
$ws.Cells.Item(9, 7).Value = 2.47329939836996
$ws.Cells.Item(9, 8).Value = 1.99947364459203
$ws.Cells.Item(9, 10).Value = 0.2616184758400806
$ws.Cells.Item(9, 11).Value = 0.5048185409211499
$ws.Cells.Item(9, 12).Value = 0.2885221836083858
$ws.Cells.Item(9, 13).Value = 0.2676901259333775

$ws.Cells.Item(10, 2).Value = 1.023849240943576
$ws.Cells.Item(10, 3).Value = 0.1154012388522574
$ws.Cells.Item(10, 4).Value = 0.06711482742188934
$ws.Cells.Item(10, 6).Value = 3.23332507531299
$ws.Cells.Item(10, 7).Value = 2.493122653842988
$ws.Cells.Item(10, 8).Value = 2.001358156822192
$ws.Cells.Item(10, 10).Value = 0.2605054587836939
$ws.Cells.Item(10, 11).Value = 0.5514538885432501
$ws.Cells.Item(10, 12).Value = 0.2912557662025108
$ws.Cells.Item(10, 13).Value = 0.2763225260437068

$ws.Cells.Item(11, 2).Value = 1.047070331969167
$ws.Cells.Item(11, 3).Value = 0.1165417030788731
$ws.Cells.Item(11, 4).Value = 0.06742382716736728
$ws.Cells.Item(11, 6).Value = 3.242195834615899
$ws.Cells.Item(11, 7).Value = 2.503095257227642
$ws.Cells.Item(11, 8).Value = 2.002844082862538
$ws.Cells.Item(11, 10).Value = 0.2600505227899959
$ws.Cells.Item(11, 11).Value = 0.5729459560070325
$ws.Cells.Item(11, 12).Value = 0.2926139002693944
$ws.Cells.Item(11, 13).Value = 0.2803846550611766

$ws.Cells.Item(12, 2).Value = 1.05593946510001
$ws.Cells.Item(12, 3).Value = 0.1169700279813739
$ws.Cells.Item(12, 4).Value = 0.06754445283610977
$ws.Cells.Item(12, 6).Value = 3.245721591507831
$ws.Cells.Item(12, 7).Value = 2.507008986465223
$ws.Cells.Item(12, 8).Value = 2.003497063917024
$ws.Cells.Item(12, 10).Value = 0.2598856193219845
$ws.Cells.Item(12, 11).Value = 0.5811241661692748
$ws.Cells.Item(12, 12).Value = 0.2931446013589891
$ws.Cells.Item(12, 13).Value = 0.2819422375971286

$ws.Cells.Item(13, 2).Value = 1.054025973625215
$ws.Cells.Item(13, 3).Value = 0.1168779379621654
$ws.Cells.Item(13, 4).Value = 0.06751831370416994
$ws.Cells.Item(13, 6).Value = 3.24495484988212
$ws.Cells.Item(13, 7).Value = 2.506159987363446
$ws.Cells.Item(13, 8).Value = 2.003352418635075
$ws.Cells.Item(13, 10).Value = 0.2599208066900971
$ws.Cells.Item(13, 11).Value = 0.5793610838319978
$ws.Cells.Item(13, 12).Value = 0.2930295769647628
$ws.Cells.Item(13, 13).Value = 0.281605925656514

$ws.Cells.Item(14, 2).Value = 1.047798483957166
$ws.Cells.Item(14, 3).Value = 0.1165770125799028
$ws.Cells.Item(14, 4).Value = 0.06743367886205931
$ws.Cells.Item(14, 6).Value = 3.242482562877058
$ws.Cells.Item(14, 7).Value = 2.50341448993737
$ws.Cells.Item(14, 8).Value = 2.002895995005161
$ws.Cells.Item(14, 10).Value = 0.260036808458679
$ws.Cells.Item(14, 11).Value = 0.573617989508449
$ws.Cells.Item(14, 12).Value = 0.2926572330174366
$ws.Cells.Item(14, 13).Value = 0.2805124112337367

$ws.Cells.Item(15, 2).Value = 1.04399382691247
$ws.Cells.Item(15, 3).Value = 0.1163922259396628
$ws.Cells.Item(15, 4).Value = 0.06738230730499595
$ws.Cells.Item(15, 6).Value = 3.240989906403854
$ws.Cells.Item(15, 7).Value = 2.501750676901736
$ws.Cells.Item(15, 8).Value = 2.002628178039885
$ws.Cells.Item(15, 10).Value = 0.2601088222819747
$ws.Cells.Item(15, 11).Value = 0.5701053317797857
$ws.Cells.Item(15, 12).Value = 0.2924312957878357
$ws.Cells.Item(15, 13).Value = 0.2798451180711368

$ws.Cells.Item(16, 2).Value = 1.022342319490974
$ws.Cells.Item(16, 3).Value = 0.115326210148389
$ws.Cells.Item(16, 4).Value = 0.06709514104579739
$ws.Cells.Item(16, 6).Value = 3.232768671042237
$ws.Cells.Item(16, 7).Value = 2.492490136451465
$ws.Cells.Item(16, 8).Value = 2.001273689041113
$ws.Cells.Item(16, 10).Value = 0.2605362223309271
$ws.Cells.Item(16, 11).Value = 0.5500548912897898
$ws.Cells.Item(16, 12).Value = 0.2911693087494882
$ws.Cells.Item(16, 13).Value = 0.2760597684905193

$ws.Cells.Item(17, 2).Value = 1.009195299474982
$ws.Cells.Item(17, 3).Value = 0.1146659170251638
$ws.Cells.Item(17, 4).Value = 0.06692544521870047
$ws.Cells.Item(17, 6).Value = 3.22802209876788
$ws.Cells.Item(17, 7).Value = 2.487053673618846
$ws.Cells.Item(17, 8).Value = 2.000603690891523
$ws.Cells.Item(17, 10).Value = 0.2608115657704353
$ws.Cells.Item(17, 11).Value = 0.5378254612734565
$ws.Cells.Item(17, 12).Value = 0.2904244215710463
$ws.Cells.Item(17, 13).Value = 0.2737721373793249

$ws.Cells.Item(18, 2).Value = 1.001683405229329
$ws.Cells.Item(18, 3).Value = 0.1142838044373846
$ws.Cells.Item(18, 4).Value = 0.06683022915224512
$ws.Cells.Item(18, 6).Value = 3.225401151788063
$ws.Cells.Item(18, 7).Value = 2.484016643755041
$ws.Cells.Item(18, 8).Value = 2.000277503654019
$ws.Cells.Item(18, 10).Value = 0.2609747733074492
$ws.Cells.Item(18, 11).Value = 0.5308175435542921
$ws.Cells.Item(18, 12).Value = 0.2900067746657982
$ws.Cells.Item(18, 13).Value = 0.2724690847271845

$ws.Cells.Item(19, 2).Value = 0.9991485943347413
$ws.Cells.Item(19, 3).Value = 0.114154027187233
$ws.Cells.Item(19, 4).Value = 0.06679840182855656
$ws.Cells.Item(19, 6).Value = 3.224532496559519
$ws.Cells.Item(19, 7).Value = 2.483003794370717
$ws.Cells.Item(19, 8).Value = 2.000177229394438
$ws.Cells.Item(19, 10).Value = 0.261030863920297
$ws.Cells.Item(19, 11).Value = 0.528449278355879
$ws.Cells.Item(19, 12).Value = 0.2898672225487147
$ws.Cells.Item(19, 13).Value = 0.2720300833294544

$ws.Cells.Item(20, 2).Value = 1.010589658333572
$ws.Cells.Item(20, 3).Value = 0.1147364473480437
$ws.Cells.Item(20, 4).Value = 0.06694326266286055
$ws.Cells.Item(20, 6).Value = 3.228516083390204
$ws.Cells.Item(20, 7).Value = 2.487623092095788
$ws.Cells.Item(20, 8).Value = 2.000668889405716
$ws.Cells.Item(20, 10).Value = 0.2607817544792113
$ws.Cells.Item(20, 11).Value = 0.5391246028326009
$ws.Cells.Item(20, 12).Value = 0.2905025996531023
$ws.Cells.Item(20, 13).Value = 0.2740143426034933

$ws.Cells.Item(21, 2).Value = 1.049625594524571
$ws.Cells.Item(21, 3).Value = 0.1166654977383743
$ws.Cells.Item(21, 4).Value = 0.06745844033245163
$ws.Cells.Item(21, 6).Value = 3.243204213452771
$ws.Cells.Item(21, 7).Value = 2.504217182184078
$ws.Cells.Item(21, 8).Value = 2.003027608002441
$ws.Cells.Item(21, 10).Value = 0.2600025360225828
$ws.Cells.Item(21, 11).Value = 0.575303802030902
$ws.Cells.Item(21, 12).Value = 0.2927661548325773
$ws.Cells.Item(21, 13).Value = 0.2808330787734477

$ws.Cells.Item(22, 2).Value = 1.075579519664586
$ws.Cells.Item(22, 3).Value = 0.117905591945842
$ws.Cells.Item(22, 4).Value = 0.06781619025593955
$ws.Cells.Item(22, 6).Value = 3.25377471460682
$ws.Cells.Item(22, 7).Value = 2.515862792956625
$ws.Cells.Item(22, 8).Value = 2.005095415598674
$ws.Cells.Item(22, 10).Value = 0.2595362269782271
$ws.Cells.Item(22, 11).Value = 0.5991798784655771
$ws.Cells.Item(22, 12).Value = 0.2943411060414149
$ws.Cells.Item(22, 13).Value = 0.2854022135849164

$ws.Cells.Item(23, 2).Value = 1.061687146134034
$ws.Cells.Item(23, 3).Value = 0.117245615511834
$ws.Cells.Item(23, 4).Value = 0.06762333665239595
$ws.Cells.Item(23, 6).Value = 3.248044247598301
$ws.Cells.Item(23, 7).Value = 2.509574071239086
$ws.Cells.Item(23, 8).Value = 2.003943670888219
$ws.Cells.Item(23, 10).Value = 0.2597811802730838
$ws.Cells.Item(23, 11).Value = 0.5864157273938133
$ws.Cells.Item(23, 12).Value = 0.2934918037389664
$ws.Cells.Item(23, 13).Value = 0.2829533015113981

$ws.Cells.Item(24, 2).Value = 1.009959123810575
$ws.Cells.Item(24, 3).Value = 0.1147045683818604
$ws.Cells.Item(24, 4).Value = 0.06693520009238796
$ws.Cells.Item(24, 6).Value = 3.228292416752709
$ws.Cells.Item(24, 7).Value = 2.487365382430824
$ws.Cells.Item(24, 8).Value = 2.000639229345694
$ws.Cells.Item(24, 10).Value = 0.2607952168719123
$ws.Cells.Item(24, 11).Value = 0.5385371894459468
$ws.Cells.Item(24, 12).Value = 0.2904672223234144
$ws.Cells.Item(24, 13).Value = 0.2739048038090317

$ws.Cells.Item(25, 2).Value = 0.9560722659437602
$ws.Cells.Item(25, 3).Value = 0.1118745157368863
$ws.Cells.Item(25, 4).Value = 0.0662843771749948
$ws.Cells.Item(25, 6).Value = 3.211160594558393
$ws.Cells.Item(25, 7).Value = 2.466845818143213
$ws.Cells.Item(25, 8).Value = 1.999337729269513
$ws.Cells.Item(25, 10).Value = 0.2620740310820153
$ws.Cells.Item(25, 11).Value = 0.4878921882372538
$ws.Cells.Item(25, 12).Value = 0.2876176330240554
$ws.Cells.Item(25, 13).Value = 0.2646316800467652
